# "Added AL Project to List"
# Insert a new row at the very top of the game list and populate it with
# the new "arcadeLightController" entry. Everything below shifts down by
# one row (values, the generated "Gx=" formula column, and the cmd-copy
# notes that live in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing rows down by one to make room for the new game.
$ws.Rows.Item(1).Insert()

# New game name goes in column B; column A re-derives the "G<n>=...exe"
# helper string from it via the existing formula pattern used by the rest
# of the sheet.
$ws.Range("B1").Value = "arcadeLightController"
$ws.Range("A1").Formula = "=""G""&ROW()-1&""=""""""&INDIRECT(""B"" & ROW())&"".exe"""""""

# Re-select the full column like Excel does after an "Insert Sheet Rows"
# operation performed from the row header.
$null = $ws.Range("A:A").Select()
